# "Actualización 10 de Mayo"
# Updates the 3 "Estadisticos" sheets (1P / 2P / Final) with revised
# Blancos/Aprobados/Por_Apro/Promedio figures, and fills in the
# "Rescatables" sheet with the list of rescatable students.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Estadisticos 1P (sheet1) — columns D (Blancos), F (Aprobados),
#    G (Por_Apro), H (Promedio) change for rows 2-6.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 6
$ws1.Range("F2").Value = 30
$ws1.Range("G2").Value = 83.33
$ws1.Range("H2").Value = 6.6

$ws1.Range("D3").Value = 5
$ws1.Range("F3").Value = 26
$ws1.Range("G3").Value = 83.87
$ws1.Range("H3").Value = 7.5

$ws1.Range("D4").Value = 4
$ws1.Range("F4").Value = 17
$ws1.Range("G4").Value = 80.95
$ws1.Range("H4").Value = 7.5

$ws1.Range("D5").Value = 7
$ws1.Range("F5").Value = 33
$ws1.Range("G5").Value = 82.5
$ws1.Range("H5").Value = 7.8

$ws1.Range("D6").Value = 4
$ws1.Range("F6").Value = 19
$ws1.Range("G6").Value = 82.61
$ws1.Range("H6").Value = 6.9

# ---------------------------------------------------------------------
# 2) Estadisticos 2P (sheet2) — columns D (Blancos), E (Reprobados),
#    F (Aprobados), G (Por_Apro) change for rows 2-6, and a new
#    H (Promedio) column is added.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 6
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 30
$ws2.Range("G2").Value = 83.33
$ws2.Range("H2").Value = 6.6

$ws2.Range("D3").Value = 5
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 26
$ws2.Range("G3").Value = 83.87
$ws2.Range("H3").Value = 7.5

$ws2.Range("D4").Value = 4
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 17
$ws2.Range("G4").Value = 80.95
$ws2.Range("H4").Value = 7.5

$ws2.Range("D5").Value = 7
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 33
$ws2.Range("G5").Value = 82.5
$ws2.Range("H5").Value = 7.8

$ws2.Range("D6").Value = 4
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 19
$ws2.Range("G6").Value = 82.61
$ws2.Range("H6").Value = 6.9

# ---------------------------------------------------------------------
# 3) Estadisticos Final (sheet3) — columns D (Blancos), F (Aprobados),
#    G (Por_Apro), H (Promedio) change for rows 2-6.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 6
$ws3.Range("F2").Value = 30
$ws3.Range("G2").Value = 83.33
$ws3.Range("H2").Value = 6.9

$ws3.Range("D3").Value = 5
$ws3.Range("F3").Value = 26
$ws3.Range("G3").Value = 83.87
$ws3.Range("H3").Value = 7.7

$ws3.Range("D4").Value = 4
$ws3.Range("F4").Value = 17
$ws3.Range("G4").Value = 80.95
$ws3.Range("H4").Value = 7.9

$ws3.Range("D5").Value = 7
$ws3.Range("F5").Value = 33
$ws3.Range("G5").Value = 82.5
$ws3.Range("H5").Value = 7.8

$ws3.Range("D6").Value = 4
$ws3.Range("F6").Value = 19
$ws3.Range("G6").Value = 82.61
$ws3.Range("H6").Value = 7.1

# ---------------------------------------------------------------------
# 4) Rescatables (sheet4) — append the 19 rescatable students. Entered
#    column-by-column (Mat, then Paterno, then Materno, ...) the same
#    way the original spreadsheet was filled in (pasted columns).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$colA_Mat = @(20330051920061, 18330051920224, 18330051920226, 18330051920220, 17330051420363, 18330051920228, 18330051920240, 18330051920244, 18330051920246, 18330061460060, 18330051920250, 17330051920302, 18330051920256, 18330051920274, 18330051920305, 18330051920283, 18330051920286, 18330051920289, 18330051920292)
$colB_Paterno = @("ROBLES", "ARIAS", "ARIAS", "ANTONIO", "ANTONIO", "BAEZ", "CERONIO", "CRESCENCIO", "GIRON", "GONZALEZ", "LOPEZ", "MARTINEZ", "MAYAHUA", "ROMAN", "TELE", "TEHUINTLE", "URBINA", "VASQUEZ", "YOPIHUA")
$colC_Materno = @("IXMATLAHUA", "GONZALEZ", "MARCELINO", "HERNANDEZ", "SANCHEZ", "MATEOS", "GARCIA", "MONTES", "CUEVAS", "RAMIREZ", "MARTINEZ", "GONZALEZ", "TEMOXTLE", "ISIDRO", "HUERTA", "MAYAHUA", "TREJO", "TORRES", "IXMATLAHUA")
$colD_Nombres = @("ALAN URIEL", "ADALI", "MIRIAM AMERICA", "VICTOR MANUEL", "MARI CARMEN", "JESUS ALFONSO", "MARIA FERNANDA", "SANDRA", "ISABEL", "EMIR ANDRES", "ELIZABETH", "LUZ GUADALUPE", "MARIA FERNANDA", "PAMELA DENISSE", "SAHARA DE JESUS", "CARLOS", "ARACELI", "KARINA", "AMAYRANI")
$colE_NombreLargo = @("INGLÉS II", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA", "TEMAS DE BIOLOGÍA CONTEMPORÁNEA")
$colF_Grupo = @("2AEV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV", "6ALCV")
$colG_NC = @(2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2, 2)

for ($i = 0; $i -lt $colA_Mat.Count; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $colA_Mat[$i]
}
for ($i = 0; $i -lt $colB_Paterno.Count; $i++) {
    $ws4.Cells.Item($i + 2, 2).Value = $colB_Paterno[$i]
}
for ($i = 0; $i -lt $colC_Materno.Count; $i++) {
    $ws4.Cells.Item($i + 2, 3).Value = $colC_Materno[$i]
}
for ($i = 0; $i -lt $colD_Nombres.Count; $i++) {
    $ws4.Cells.Item($i + 2, 4).Value = $colD_Nombres[$i]
}
for ($i = 0; $i -lt $colE_NombreLargo.Count; $i++) {
    $ws4.Cells.Item($i + 2, 5).Value = $colE_NombreLargo[$i]
}
for ($i = 0; $i -lt $colF_Grupo.Count; $i++) {
    $ws4.Cells.Item($i + 2, 6).Value = $colF_Grupo[$i]
}
for ($i = 0; $i -lt $colG_NC.Count; $i++) {
    $ws4.Cells.Item($i + 2, 7).Value = $colG_NC[$i]
}
